$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 to 39/51: price (D) and volume% (E) updates ---
$ws.Range("D2").Value = "29.592.34"
$ws.Range("E2").Value = "  +0.95%  "

$ws.Range("D3").Value = "1.987.40"
$ws.Range("E3").Value = "  +4.52%  "

$ws.Range("D4").Value = "1.004"
$ws.Range("E4").Value = "  +0.40%  "

$ws.Range("D5").Value = "327.76"
$ws.Range("E5").Value = "  +0.52%  "

$ws.Range("D6").Value = "1.004"
$ws.Range("E6").Value = "  +0.40%  "

$ws.Range("D7").Value = "0.4663"
$ws.Range("E7").Value = "  +0.45%  "

$ws.Range("D8").Value = "0.3925"
$ws.Range("E8").Value = "  +0.01%  "

$ws.Range("D9").Value = "46.12"
$ws.Range("E9").Value = "  -1.31%  "

$ws.Range("D10").Value = "0.07929"
$ws.Range("E10").Value = "  +0.45%  "

$ws.Range("D12").Value = "22.89"
$ws.Range("E12").Value = "  +3.72%  "

$ws.Range("D13").Value = "1.985.93"
$ws.Range("E13").Value = "  +5.36%  "

$ws.Range("E14").Value = "  +1.50%  "

$ws.Range("D15").Value = "5.852"
$ws.Range("E15").Value = "  +1.60%  "

$ws.Range("D16").Value = "0.07104"
$ws.Range("E16").Value = "  +1.54%  "

$ws.Range("D17").Value = "87.71"
$ws.Range("E17").Value = "  -0.88%  "

$ws.Range("E18").Value = "  +0.36%  "

$ws.Range("D19").Value = "0.000009951"

$ws.Range("E20").Value = "  +1.46%  "

$ws.Range("E21").Value = "  +0.22%  "

$ws.Range("D22").Value = "29.615.55"
$ws.Range("E22").Value = "  +1.10%  "

$ws.Range("D23").Value = "5.545"
$ws.Range("E23").Value = "  +4.43%  "

$ws.Range("D24").Value = "11.19"

$ws.Range("D25").Value = "2.226.52"
$ws.Range("E25").Value = "  +5.15%  "

$ws.Range("D26").Value = "2.107"
$ws.Range("E26").Value = "  +0.49%  "

$ws.Range("D27").Value = "159.05"
$ws.Range("E27").Value = "  +1.71%  "

$ws.Range("D28").Value = "19.58"
$ws.Range("E28").Value = "  +0.48%  "

$ws.Range("D29").Value = "5.813"
$ws.Range("E29").Value = "  -3.56%  "

$ws.Range("D30").Value = "119.59"
$ws.Range("E30").Value = "  +0.73%  "

$ws.Range("D31").Value = "1.897"
$ws.Range("E31").Value = "  -1.47%  "

$ws.Range("D32").Value = "0.09422"
$ws.Range("E32").Value = "  +0.50%  "

$ws.Range("D33").Value = "0.8953"
$ws.Range("E33").Value = "  -1.39%  "

$ws.Range("D34").Value = "5.230"
$ws.Range("E34").Value = "  -1.30%  "

$ws.Range("D35").Value = "1.334"
$ws.Range("E35").Value = "  +0.29%  "

$ws.Range("D36").Value = "3.195"
$ws.Range("E36").Value = "  -1.00%  "

$ws.Range("D37").Value = "0.05805"
$ws.Range("E37").Value = "  +0.03%  "

$ws.Range("E38").Value = "  -0.96%  "

$ws.Range("D39").Value = "0.02100"
$ws.Range("E39").Value = "  +0.22%  "

$ws.Range("D42").Value = "0.5740"
$ws.Range("E42").Value = "  +0.38%  "

$ws.Range("D43").Value = "0.1804"
$ws.Range("E43").Value = "  +0.88%  "

$ws.Range("E44").Value = "  -1.02%  "

$ws.Range("D45").Value = "2.783"
$ws.Range("E45").Value = "  +7.94%  "

$ws.Range("D46").Value = "11.89"
$ws.Range("E46").Value = "  -0.46%  "

$ws.Range("D47").Value = "0.5361"
$ws.Range("E47").Value = "  +0.00%  "

$ws.Range("D48").Value = "2.169"
$ws.Range("E48").Value = "  -2.17%  "

$ws.Range("D49").Value = "0.06943"
$ws.Range("E49").Value = "  -1.59%  "

$ws.Range("D50").Value = "114.21"
$ws.Range("E50").Value = "  +0.96%  "

$ws.Range("D51").Value = "1.825"
$ws.Range("E51").Value = "  -2.16%  "

# --- Row 40 and 41: swap PEPE and FraxShare, with updated values ---
$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D40").Value = "7.846"
$ws.Range("E40").Value = "  +0.62%  "

$ws.Range("B41").Value = "PEPE"
$ws.Range("C41").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D41").Value = "0.000003251"
$ws.Range("E41").Value = "  +46.58%  "
